$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new dSF (column F) value, per repulled/pushed data
$updates = @{
    2  = -10
    6  = 0
    8  = 12
    10 = -4
    11 = -3
    12 = -1
    13 = -9
    15 = -6
    17 = -3
    19 = -1
    20 = 5
    22 = -1
    23 = -6
    24 = 2
    25 = 4
    26 = -4
    27 = -3
    29 = 3
    31 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
